$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing existing rows 20-32 down to 21-33.
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with this week's record.
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(20, 4).Value = 44942
$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100112003
$ws.Cells.Item(20, 7).Value = "Ajo"
$ws.Cells.Item(20, 8).Value = "Chino"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 1000
$ws.Cells.Item(20, 11).Value = 14000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 14500
$ws.Cells.Item(20, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(20, 15).Value = "China"
$ws.Cells.Item(20, 16).Value = 1450
$ws.Cells.Item(20, 17).Value = 10
$ws.Cells.Item(20, 18).Value = "Hortaliza"
